# Update consolidated_reports.xlsx "Reports" sheet:
# - Update file_path values to the new D:\hackathon-evaluation\... prefix
# - Refresh evaluation summaries/feedback text (model re-run, e.g. gpt-4o)
# - Reorder/replace team rows 3-6 and append 5 new team rows (7-11)
# - Resulting used range becomes A1:Q11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Coding Pirates'
$ws.Range("B2").Value = 'D:\hackathon-evaluation\hackathon-evaluation\project_context\ppt\2025CodingPirates - YASH KASAUDHAN.pdf'
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 37
$ws.Range("K2").Value = 61.5
$ws.Range("L2").Value = 'The project by Coding Pirates aims to digitize and showcase the monasteries of Sikkim for tourism and cultural preservation. The deck provides a comprehensive overview of the problem and the proposed solution, including a detailed technical architecture. However, there are gaps in baseline data, cost estimates, and security considerations. The project demonstrates strong potential impact and a clear adoption path, but lacks detailed metrics and evaluation plans.'
$ws.Range("M2").Value = 'The workflow involves creating a digital platform for virtual tours and cultural preservation of Sikkim''s monasteries. It includes a frontend for user interaction, a backend for data management, and AI components for enhanced search capabilities.'
$ws.Range("N2").Value = '1. The project addresses a significant cultural and tourism challenge with a well-defined problem statement. 2. The technical architecture is robust and leverages mature technologies, which enhances feasibility.'
$ws.Range("O2").Value = '1. The deck lacks detailed baseline data and cost estimates, which are crucial for assessing the project''s viability. 2. Security considerations are minimally addressed, posing potential risks.'
$ws.Range("P2").Value = '1. The use of React, Node.js, and PostGIS is appropriate for the project''s requirements, ensuring scalability and performance. 2. The AI components for image recognition and NLP search are well-integrated but require further validation.'
$ws.Range("Q2").Value = '1. Include detailed metrics and evaluation plans to measure the project''s success effectively. 2. Enhance the security framework to protect user data and ensure compliance with data protection laws.'

# Row 3
$ws.Range("A3").Value = 'ALT_F4'
$ws.Range("B3").Value = 'D:\hackathon-evaluation\hackathon-evaluation\project_context\ppt\ALT_f4 - VAIBHAV KUMAR.pdf'
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 32
$ws.Range("K3").Value = 54
$ws.Range("L3").Value = 'The team proposes an AI-driven chatbot for public health awareness, focusing on disease prevention and multilingual support. The solution leverages NLP and ML models, integrating trusted medical data sources. While the idea is innovative, the presentation lacks detailed metrics, baselines, and a comprehensive evaluation plan. The architecture and scalability are addressed at a high level, but specifics on latency, cost, and security are minimal.'
$ws.Range("M3").Value = 'No diagrams were provided, so the workflow is inferred from the text. The process involves data collection, preprocessing, NLP model training, chatbot integration, and deployment. The chatbot is designed to be multilingual and continuously learns from user feedback.'
$ws.Range("N3").Value = '1. The integration of AI with vernacular language support is a strong point, addressing a significant need in diverse linguistic regions. 2. The use of trusted medical data sources enhances the credibility of the chatbot.'
$ws.Range("O3").Value = '1. The presentation lacks diagrams, which makes it difficult to visualize the architecture and workflow. 2. There is insufficient detail on the evaluation plan and metrics to measure success.'
$ws.Range("P3").Value = '1. The choice of technologies like TensorFlow/PyTorch and cloud deployment is appropriate, but more details on the specific architecture would be beneficial. 2. The security and privacy aspects are not thoroughly addressed, which is critical for handling sensitive health data.'
$ws.Range("Q3").Value = '1. Include detailed diagrams to illustrate the architecture and workflow. 2. Provide a comprehensive evaluation plan with clear metrics and baselines. 3. Address security and privacy concerns more thoroughly, possibly with a dedicated section on compliance with health data regulations.'

# Row 4
$ws.Range("A4").Value = 'Kairos'
$ws.Range("B4").Value = 'D:\hackathon-evaluation\hackathon-evaluation\project_context\ppt\AgriNiti - TEENA gla.pdf'
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 6
$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 36
$ws.Range("K4").Value = 60
$ws.Range("L4").Value = 'The team Kairos presents a solution aimed at empowering small and marginal farmers through an AI-driven crop advisory system. The solution integrates soil, climate, and market data to provide crop recommendations, pest and disease detection, and localized weather alerts. The approach leverages standard web technologies and government APIs, with a focus on accessibility through a multilingual mobile app. However, the presentation lacks detailed diagrams and metrics, and the assumptions and baselines are not well-defined.'
$ws.Range("M4").Value = 'No diagrams were provided, so the workflow could not be analyzed.'
$ws.Range("N4").Value = '1. The solution addresses a significant problem for small and marginal farmers, aiming to increase income and reduce crop loss. 2. The use of AI for crop recommendations and pest detection is innovative and has the potential for high impact.'
$ws.Range("O4").Value = '1. The presentation lacks detailed diagrams, which makes it difficult to assess the technical feasibility and architecture of the solution. 2. There is insufficient evidence of a well-defined evaluation plan or metrics to measure success.'
$ws.Range("P4").Value = '1. The integration with existing government APIs is a strong point, but the scalability and latency of the system are not addressed. 2. The use of standard web technologies is appropriate, but more details on the backend architecture are needed.'
$ws.Range("Q4").Value = '1. Include detailed diagrams to illustrate the system architecture and data flow. 2. Define clear metrics and an evaluation plan to measure the impact and success of the solution. 3. Address potential risks and provide mitigation strategies, especially concerning data privacy and security.'

# Row 5
$ws.Range("A5").Value = 'Algo Wizards'
$ws.Range("B5").Value = 'D:\hackathon-evaluation\hackathon-evaluation\project_context\ppt\Algo wizards - LAXMI gla.pdf'
$ws.Range("D5").Value = 6
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 6
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 35
$ws.Range("K5").Value = 58.5
$ws.Range("L5").Value = 'Algo Wizards present an AI-based Farmer Query Support and Advisory System aimed at addressing the challenges faced by farmers in accessing real-time, accurate advice. The solution leverages multilingual NLP and voice-based interaction to break language barriers and provide instant, expert-level guidance. While the problem is well-framed and the solution innovative, the deck lacks detailed baselines, metrics, and a comprehensive evaluation plan. The architecture is modular, allowing scalability, but lacks detailed latency and cost estimates. Privacy, compliance, and security considerations are minimally addressed.'
$ws.Range("M5").Value = 'The workflow involves farmers submitting queries via voice, text, or image, which are then preprocessed and analyzed by an AI engine. The system generates advisory outputs in local languages, providing guidance on crop diseases, fertilizer recommendations, weather-based irrigation, and government schemes.'
$ws.Range("N5").Value = '1. The solution effectively addresses language barriers and accessibility issues for non-literate users through multilingual NLP and voice-based interaction. 2. The modular architecture allows for easy scaling across different regions and crop types.'
$ws.Range("O5").Value = '1. The deck lacks detailed baselines and metrics to evaluate the effectiveness of the solution. 2. Privacy, compliance, and security considerations are not adequately addressed.'
$ws.Range("P5").Value = '1. The technical stack is well-defined, but the absence of latency and cost estimates raises concerns about the feasibility of large-scale deployment. 2. The use of APIs for weather, government schemes, and translation is a strong point, but the integration details are not provided.'
$ws.Range("Q5").Value = '1. Include detailed baselines and metrics to measure the impact and effectiveness of the solution. 2. Address privacy, compliance, and security concerns to ensure user trust and data protection. 3. Provide latency and cost estimates to assess the feasibility of scaling the solution.'

# Row 6
$ws.Range("A6").Value = 'AlgoYoddhas'
$ws.Range("B6").Value = 'D:\hackathon-evaluation\hackathon-evaluation\project_context\ppt\AlgoYoddhas - ADWAIT PATEL.pdf'
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 28
$ws.Range("K6").Value = 46
$ws.Range("L6").Value = 'The AlgoYoddhas team proposes a mobile and web platform for crowdsourced ocean hazard reporting, integrating social media analytics. The solution includes geotagged reports, a real-time dashboard, and an NLP engine for sentiment analysis. While the proposal outlines a feasible approach using open-source GIS and NLP tools, it lacks detailed diagrams and metrics. The team identifies challenges such as data verification and limited connectivity, offering solutions like role-based access and AI filters. However, the absence of diagrams and specific metrics limits the evaluation of their technical approach and scalability.'
$ws.Range("M6").Value = 'No diagrams were provided, so the workflow analysis is based solely on the text description. The proposed workflow involves users submitting geotagged reports, which are then analyzed by an NLP engine and displayed on a real-time dashboard. The platform supports multilingual and offline capabilities, aiming to enhance situational awareness and response times.'
$ws.Range("N6").Value = '1. The platform addresses a critical need for real-time hazard reporting and social media analytics in disaster management. 2. The use of open-source GIS and NLP tools is a practical approach to ensure feasibility and cost-effectiveness.'
$ws.Range("O6").Value = '1. The proposal lacks detailed diagrams, which makes it difficult to fully understand the technical architecture and workflow. 2. There is insufficient information on the datasets and baselines used for the NLP engine.'
$ws.Range("P6").Value = '1. The absence of latency and cost estimates raises concerns about the platform''s scalability and performance in real-world scenarios. 2. The security and privacy compliance aspects are not adequately addressed, which are critical for handling sensitive data.'
$ws.Range("Q6").Value = '1. Include detailed architecture diagrams to better illustrate the technical workflow and system components. 2. Provide specific metrics and KPIs to evaluate the platform''s effectiveness and impact. 3. Address privacy and security concerns with clear compliance strategies.'

# Row 7
$ws.Range("A7").Value = 'Alt-Era'
$ws.Range("B7").Value = 'D:\hackathon-evaluation\hackathon-evaluation\project_context\ppt\Alt-Era - KRISH PATHAK.pdf'
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 7
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 6
$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 38
$ws.Range("K7").Value = 64
$ws.Range("L7").Value = 'Alt-Era presents a solution for early detection of water-borne diseases using a mobile app and SMS-based platform integrated with IoT and AI/ML technologies. The solution is innovative in its combination of health and water data for outbreak prediction, and it is designed to work in low-bandwidth environments with multilingual support. However, the deck lacks detailed baselines, datasets, and privacy/security considerations. The architecture is well-defined, but cost and latency estimates are missing.'
$ws.Range("M7").Value = 'No diagrams were provided, so the workflow analysis is based solely on the text description. The proposed workflow involves data collection through mobile apps, SMS, and IoT sensors, followed by AI/ML-driven data processing for outbreak detection and notification to officials and communities.'
$ws.Range("N7").Value = '1. The solution is innovative, combining health and water data with AI/ML for early outbreak prediction. 2. The platform is designed to work in low-bandwidth and offline environments, which is crucial for rural areas. 3. Multilingual support enhances accessibility for tribal communities.'
$ws.Range("O7").Value = '1. The deck lacks diagrams, which makes it difficult to visualize the architecture and workflow. 2. There is insufficient detail on baselines and datasets used for AI/ML models. 3. Privacy and security considerations are not addressed in the deck.'
$ws.Range("P7").Value = '1. The use of low-cost IoT sensors and existing ASHA workforce is a practical approach to data collection. 2. The backend technologies (Django/Node.js) and AI/ML frameworks (TensorFlow, Scikit-learn) are appropriate for the solution.'
$ws.Range("Q7").Value = '1. Include diagrams to illustrate the architecture and data flow. 2. Provide more details on the datasets and baselines for AI/ML models. 3. Address privacy and security concerns, especially given the sensitive nature of health data.'

# Row 8
$ws.Range("A8").Value = 'Angaari Paltan'
$ws.Range("B8").Value = 'D:\hackathon-evaluation\hackathon-evaluation\project_context\ppt\Angaari Paltan - ISHITA GOYAL.pdf'
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 7
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 8
$ws.Range("J8").Value = 38
$ws.Range("K8").Value = 63.5
$ws.Range("L8").Value = 'The team ''Angaari Paltan'' presented an AI-driven internship recommendation engine aimed at matching user profiles with verified internships. The solution leverages NLP for profile analysis and a hybrid model for matching. However, the presentation lacks diagrams, making it difficult to fully assess the technical architecture and workflow. The problem framing and datasets are adequately covered, but there are significant gaps in assumptions, baselines, and metrics. The deployment plan and adoption path are mentioned but not detailed. Overall, the project shows potential but requires more concrete evidence and detailed planning.'
$ws.Range("M8").Value = 'No diagrams were provided, so the workflow analysis is based solely on the text description. The workflow involves data collection from verified sources, user profiling using NLP, and a recommendation engine that uses a hybrid model for matching. The deployment includes a web dashboard for ministries and organizations.'
$ws.Range("N8").Value = '1. The project addresses a relevant problem in the education sector by facilitating access to internships. 2. The use of NLP for profile analysis is a strong technical choice that can enhance the personalization of recommendations.'
$ws.Range("O8").Value = '1. The lack of diagrams makes it difficult to assess the technical architecture and workflow. 2. There is insufficient detail on the assumptions and baselines, which are critical for evaluating the feasibility of the solution.'
$ws.Range("P8").Value = '1. The hybrid model combining collaborative filtering and skill matching is a promising approach, but more details on its implementation and evaluation are needed. 2. The data quality challenges are acknowledged, but the proposed strategies for verification and feedback loops require further elaboration.'
$ws.Range("Q8").Value = '1. Include detailed diagrams to illustrate the architecture and workflow, which will strengthen the technical evaluation. 2. Provide more information on the metrics and evaluation plan to demonstrate how the solution''s effectiveness will be measured. 3. Clarify the assumptions and baselines to provide a clearer context for the solution''s development and deployment.'

# Row 9
$ws.Range("A9").Value = 'Hacktronics'
$ws.Range("B9").Value = 'D:\hackathon-evaluation\hackathon-evaluation\project_context\ppt\Animal Type Classification - DEVANG SHUKLA.pdf'
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 8
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 6
$ws.Range("I9").Value = 8
$ws.Range("J9").Value = 39
$ws.Range("K9").Value = 65.5
$ws.Range("L9").Value = 'Hacktronics presents an AI-driven solution for animal type classification in rural agriculture, leveraging low-cost hardware and AI frameworks. The project is well-framed with a clear problem statement and a strong architectural plan, but lacks detailed baselines, metrics, and privacy considerations. The diagrams support the data flow and market potential, but the text lacks comprehensive evidence for some criteria.'
$ws.Range("M9").Value = 'The workflow involves capturing data using a camera, processing it on-device, and transmitting it to a mobile device and cloud, supported by market sizing data for deployment scenarios.'
$ws.Range("N9").Value = '1. The project addresses a significant problem in rural agriculture with a clear AI-driven solution. 2. The use of low-cost hardware and AI frameworks is well-integrated into the solution, making it feasible for rural deployment.'
$ws.Range("O9").Value = '1. The proposal lacks detailed baselines and metrics for evaluating the system''s performance. 2. Privacy and security considerations are minimally addressed, which could be a concern for data handling.'
$ws.Range("P9").Value = '1. The architecture is well-documented, but the latency and cost estimates are vague. 2. The dataset requirements are mentioned, but specifics on data collection and diversity are lacking.'
$ws.Range("Q9").Value = '1. Provide more detailed metrics and evaluation plans to strengthen the proposal. 2. Address privacy and security concerns more thoroughly to ensure data protection. 3. Include more detailed deployment and adoption strategies to enhance the project''s feasibility.'

# Row 10
$ws.Range("A10").Value = 'BenzeneCoder'
$ws.Range("B10").Value = 'D:\hackathon-evaluation\hackathon-evaluation\project_context\ppt\BenzeneCoder - SARTHAK TIWARI.pdf'
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 7
$ws.Range("G10").Value = 4
$ws.Range("H10").Value = 6
$ws.Range("I10").Value = 8
$ws.Range("J10").Value = 38
$ws.Range("K10").Value = 64
$ws.Range("L10").Value = 'The team BenzeneCoder presents a digital mental health platform for students, focusing on accessibility and stigma reduction. The solution is technically feasible using open-source tools but lacks detailed metrics, baselines, and deployment plans. Privacy and security are acknowledged but not deeply addressed.'
$ws.Range("M10").Value = 'The workflow involves developing a web-based platform with core features like AI chatbot, booking system, and resource hub. The process includes defining features, setting up backend and frontend, integrating components, and testing.'
$ws.Range("N10").Value = '1. The project addresses a significant issue in student mental health with a proactive approach. 2. The use of open-source tools makes the solution accessible and cost-effective.'
$ws.Range("O10").Value = '1. The presentation lacks detailed diagrams to support the technical architecture and workflow. 2. Metrics and evaluation plans are not clearly defined, making it difficult to assess the project''s success.'
$ws.Range("P10").Value = '1. The integration of AI chatbot with Rasa is a strong point, but the limitations and potential inaccuracies of AI need more robust mitigation strategies. 2. The security measures for handling sensitive data are mentioned but not detailed.'
$ws.Range("Q10").Value = '1. Include detailed diagrams of the system architecture and data flow to strengthen the technical presentation. 2. Develop a clear evaluation plan with specific metrics to measure the impact and success of the platform. 3. Enhance the privacy and security sections with more specific strategies and technologies.'

# Row 11
$ws.Range("A11").Value = 'Binary Brains'
$ws.Range("B11").Value = 'D:\hackathon-evaluation\hackathon-evaluation\project_context\ppt\Binary Brains - Milan Sharma.pdf'
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 6
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 7
$ws.Range("J11").Value = 33
$ws.Range("K11").Value = 55.5
$ws.Range("L11").Value = 'Binary Brains presents a digital mental health platform for students, leveraging AI for stress detection and personalized support. The proposal is strong in problem framing and architecture, but lacks detailed baselines and dataset information. Privacy and scalability are well-addressed, but metrics and evaluation plans need more clarity.'
$ws.Range("M11").Value = 'No diagrams or workflows were provided in the deck. The images were primarily logos and decorative elements.'
$ws.Range("N11").Value = '1. The problem framing is well-articulated, focusing on a critical issue in student mental health. 2. The proposed solution is innovative, combining AI with human counseling to provide comprehensive support.'
$ws.Range("O11").Value = '1. The deck lacks detailed information on datasets and baselines, which are crucial for evaluating the AI model''s effectiveness. 2. Metrics and evaluation plans are not clearly defined, making it difficult to assess the project''s success criteria.'
$ws.Range("P11").Value = '1. The tech stack is well-chosen for scalability and real-time interaction, but more details on latency and cost estimates would strengthen the proposal. 2. Security measures are mentioned, but specifics on how data privacy will be maintained are needed.'
$ws.Range("Q11").Value = '1. Include detailed baselines and dataset information to strengthen the evidence of feasibility. 2. Define clear metrics and evaluation plans to track the project''s impact and success. 3. Provide a more detailed deployment plan, including potential challenges and mitigation strategies.'
